# Code Version 29052020 0234
#
# The TSYS_PJA_InputFile_Pattern column (column O) holds regex patterns used
# to match incoming "TransactionMonetary_*" file names. The timestamp portion
# of those patterns is widened from 8 digits to 14 digits
# (\d\d\d\d\d\d\d\d -> \d\d\d\d\d\d\d\d\d\d\d\d\d\d), i.e. the date-only
# "yyyyMMdd" stamp becomes a full "yyyyMMddHHmmss" stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new value for column O ("TSYS_PJA_InputFile_Pattern").
$updates = @{
    2  = '^TransactionMonetary_maynilad_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    4  = '^TransactionMonetary_veco_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    5  = '^TransactionMonetary_branch_payment_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    7  = '^TransactionMonetary_auto_debit_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    11 = '^TransactionMonetary_cash_advance_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    12 = '^TransactionMonetary_branch_payment_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    13 = '^TransactionMonetary_branch_payment_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    17 = '^TransactionMonetary_pldt_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    18 = '^TransactionMonetary_branch_payment_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    19 = '^TransactionMonetary_bpi_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    21 = '^TransactionMonetary_meralco_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    23 = '^TransactionMonetary_smart_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    24 = '^TransactionMonetary_sun_cellular_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    25 = '^TransactionMonetary_sm_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    27 = '^TransactionMonetary_sky_cable_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    28 = '^TransactionMonetary_globe_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
    29 = '^TransactionMonetary_bayantel_\d\d\d\d\d\d\d\d\d\d\d\d\d\d.dat$'
}

foreach ($r in $updates.Keys) {
    $ws.Cells.Item($r, 15).Value2 = $updates[$r]   # Column O = 15 = TSYS_PJA_InputFile_Pattern
}

# The author's final selection, per the saved workbook view state.
$ws.Range("C1").Select()
